$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.435.30"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.688.29"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "680.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "4.309.48"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "3.682.21"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "69.402.56"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.10"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "3.833.43"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.91"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.676.58"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.162"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  +2.47%  "
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "169.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.939"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.57"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000277"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.57%  "
